$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; this shifts the existing rows 17-35 down to 18-36
# and copies formatting (e.g. the date number format on column D) from the row below.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new record's data.
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = "2022-01-25"
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100103
$ws.Range("H17").Value = "Frutos de hueso (carozo)"
$ws.Range("I17").Value = 100103002
$ws.Range("J17").Value = "Ciruela"
$ws.Range("K17").Value = "Black Amber"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 9500
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 9750
$ws.Range("Q17").Value = "$/bandeja 18 kilos granel"
$ws.Range("R17").Value = "Región del Maule"
$ws.Range("S17").Value = 542
$ws.Range("T17").Value = 18
